# The post previously stored in row 668 ("「睡眠の素晴らしさに匹敵するものは無い！」")
# was removed from the source data. Delete that row entirely and let Excel
# shift all the subsequent rows (669-843) up by one, which naturally
# shrinks the used range from A1:C843 to A1:C842.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(668).Delete()
